# Update Betfair back/lay odds for 2025-11-20 games (rows 2-4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1.87
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3.5
$ws.Range("L2").Value = 1.37
$ws.Range("N2").Value = 3.2
$ws.Range("P2").Value = 1.77
$ws.Range("Q2").Value = 2.06
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 3.9
$ws.Range("T2").Value = 1.94
$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 2
$ws.Range("X2").Value = 15
$ws.Range("Z2").Value = 11.5
$ws.Range("AA2").Value = 23
$ws.Range("AB2").Value = 15
$ws.Range("AE2").Value = 23
$ws.Range("AF2").Value = 980
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 980
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 80
$ws.Range("AM2").Value = 160
$ws.Range("AN2").Value = 110
$ws.Range("AO2").Value = 16.5
$ws.Range("F3").Value = 1.54
$ws.Range("G3").Value = 1.63
$ws.Range("K3").Value = 4.9
$ws.Range("L3").Value = 1.28
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.4
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.74
$ws.Range("R3").Value = 1.47
$ws.Range("S3").Value = 2.88
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 2.02
$ws.Range("W3").Value = 2.6
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 55
$ws.Range("AA3").Value = 220
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 10.5
$ws.Range("AD3").Value = 980
$ws.Range("AE3").Value = 110
$ws.Range("AF3").Value = 980
$ws.Range("AG3").Value = 10.5
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 980
$ws.Range("AK3").Value = 980
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 130
$ws.Range("AN3").Value = 8.6
$ws.Range("F4").Value = 2.44
$ws.Range("G4").Value = 2.64
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.75
$ws.Range("J4").Value = 2.96
$ws.Range("K4").Value = 3.15
$ws.Range("N4").Value = 2.68
$ws.Range("P4").Value = 1.57
$ws.Range("Q4").Value = 2.52
$ws.Range("R4").Value = 1.2
$ws.Range("S4").Value = 5.1
$ws.Range("T4").Value = 2.04
$ws.Range("V4").Value = 1.36
$ws.Range("W4").Value = 1.61
$ws.Range("Z4").Value = 980
$ws.Range("AD4").Value = 980
$ws.Range("AF4").Value = 980
$ws.Range("AG4").Value = 980
$ws.Range("AH4").Value = 980
$ws.Range("AI4").Value = 80
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 980
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 980
